$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of profit-allocation data appended after the 2025-11-17 run.
# Force the date column to Text first so the COM layer's "smart" date
# parser doesn't turn "11/17/2025" into a date serial number (the source
# data stores dates as plain text, matching every prior row), then reset
# the style back to Normal so no stray per-cell formatting is introduced.
$ws.Range("A77").NumberFormat = "@"
$ws.Range("A77").Value = "11/17/2025"
$ws.Range("A77").Style = "Normal"

$ws.Range("B77").Value = 0.2014044958689747
$ws.Range("C77").Value = 0.7985955041310253
